# Update imputed values in the KNN result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B11").Value = 6.447
$ws.Range("A12").Value = -21.489
$ws.Range("B23").Value = 8.695
$ws.Range("D24").Value = -7.635999999999998
$ws.Range("B28").Value = 5.196000000000001
$ws.Range("A32").Value = -21.043
$ws.Range("B32").Value = 6.795
$ws.Range("B34").Value = 7.119999999999999
$ws.Range("A36").Value = -20.724
$ws.Range("A38").Value = -20.417
$ws.Range("D38").Value = -8.348000000000001
$ws.Range("B42").Value = 9.103999999999999
$ws.Range("A46").Value = -21.578
$ws.Range("D52").Value = -7.941000000000001
$ws.Range("A54").Value = -21.185
$ws.Range("B54").Value = 6.145
$ws.Range("A55").Value = -22.016
$ws.Range("A67").Value = -21.422
$ws.Range("A69").Value = -21.387
$ws.Range("A72").Value = -21.621
$ws.Range("D78").Value = -8.253
$ws.Range("D83").Value = -8.141000000000002
$ws.Range("D85").Value = -8.643000000000001
$ws.Range("D86").Value = -8.282
$ws.Range("A91").Value = -20.847
$ws.Range("D96").Value = -7.253
$ws.Range("B97").Value = 5.167
$ws.Range("A99").Value = -21.157
$ws.Range("B99").Value = 6.016000000000001
$ws.Range("B101").Value = 5.252
$ws.Range("D103").Value = -8.341999999999999
$ws.Range("A104").Value = -21.437
